$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.24711565198759189
$ws.Range("B1").Value = 0.24692851511571234
$ws.Range("A2").Value = -0.17176775076050887
$ws.Range("B2").Value = 0.17129964859223001
$ws.Range("A3").Value = -0.12159800624577777
$ws.Range("B3").Value = 0.12124860204961863
$ws.Range("A4").Value = -0.11324860214501875
$ws.Range("B4").Value = 0.11278319271417026
$ws.Range("A5").Value = -0.10978319277053217
$ws.Range("B5").Value = 0.1081962704864905
$ws.Range("A6").Value = -0.008998037416754201
$ws.Range("B6").Value = 0.0088932779540460416
$ws.Range("A7").Value = 0.0011067219086622337
$ws.Range("B7").Value = -0.0011202998511401496
$ws.Range("A8").Value = 0.011120299713931026
$ws.Range("B8").Value = -0.011131093353121901
$ws.Range("A9").Value = 0.013131093290520202
$ws.Range("B9").Value = -0.013140017285891759
$ws.Range("A10").Value = 0.015140017225132141
$ws.Range("B10").Value = -0.015139645575001026
$ws.Range("A11").Value = -0.03212773405809255
$ws.Range("B11").Value = 0.032098742104209599
$ws.Range("A12").Value = -0.028598742180213854
$ws.Range("B12").Value = 0.028377273312641282
$ws.Range("A13").Value = -0.017164761379098614
$ws.Range("B13").Value = 0.017079145135324936
$ws.Range("A14").Value = -0.0090791452610003986
$ws.Range("B14").Value = 0.0090514825881342631
$ws.Range("A15").Value = -0.0080514826507451787
$ws.Range("B15").Value = 0.0080336412915515965
$ws.Range("A16").Value = -0.0060336413642456677
$ws.Range("B16").Value = 0.0060031678911722786
$ws.Range("A17").Value = -0.0040031679652372532
$ws.Range("B17").Value = 0.0039999999072559689
$ws.Range("A18").Value = -0.069160764413009446
$ws.Range("B18").Value = 0.069052228812378047
$ws.Range("A19").Value = -0.065052228853460292
$ws.Range("B19").Value = 0.064234019721200575
$ws.Range("A20").Value = -0.06023401977918752
$ws.Range("B20").Value = 0.059999582657075834
$ws.Range("A21").Value = -0.0040057281362759767
$ws.Range("B21").Value = 0.0039999999380118112
$ws.Range("A22").Value = -0.045701642398864095
$ws.Range("B22").Value = 0.04549177192592424
$ws.Range("A23").Value = -0.040491771990414094
$ws.Range("B23").Value = 0.040097479334415098
$ws.Range("A24").Value = -0.020097479541487218
$ws.Range("B24").Value = 0.019999999790274892
$ws.Range("A25").Value = -0.097198233133507728
$ws.Range("B25").Value = 0.097077736802095416
$ws.Range("A26").Value = -0.094577736873114659
$ws.Range("B26").Value = 0.094421575534235913
$ws.Range("A27").Value = -0.091921575608995276
$ws.Range("B27").Value = 0.090990497104204593
$ws.Range("A28").Value = -0.088990497190072126
$ws.Range("B28").Value = 0.088353188064519195
$ws.Range("A29").Value = -0.08135318820895332
$ws.Range("B29").Value = 0.081168844760910908
$ws.Range("A30").Value = -0.021168845396268399
$ws.Range("B30").Value = 0.021021976700056122
$ws.Range("A31").Value = -0.014021976854902363
$ws.Range("B31").Value = 0.014000862567646877
$ws.Range("A32").Value = -0.0040008627506225025
$ws.Range("B32").Value = 0.0039999998716755414
